$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# The old "Sheet1" carries legacy column widths / merged cells / row
# formatting from the previous (First/Second/Third check-in/out) layout
# that can't be fully reset in place. Start from a brand-new sheet named
# "Attendance" instead, then drop the now unused Sheet1 / Sheet2 / Sheet3.
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "Attendance"
$wb.Worksheets.Item("Sheet1").Delete()
$wb.Worksheets.Item("Sheet2").Delete()
$wb.Worksheets.Item("Sheet3").Delete()

$ws = $wb.Worksheets.Item("Attendance")

# New simplified header row: Code | Name | Type | Time
$ws.Range("A1").Value = "Code"
$ws.Range("B1").Value = "Name"
$ws.Range("C1").Value = "Type"
$ws.Range("D1").Value = "Time"

# Thin box border around each header cell (time-log formatting).
$ws.Range("A1:D1").Borders.LineStyle = 1

$ws.Range("F2").Select()
